# change tracing strategy and save wallet labels
# Append two new rows (6 and 7) to Sheet1 with date labels + USD values.
#
# Note: assigning an ISO-looking date string ("2024-10-04") directly to a
# cell's .Value triggers Excel's normal "smart entry" and the cell becomes
# a date serial number with a date number-format style — but the existing
# rows in this sheet store their dates as plain text (shared strings) with
# no special style. To reproduce that, we build the text via a formula
# (which is never auto-converted), then Copy + PasteSpecial(xlPasteValues)
# it into the target cell — pasting a text formula result keeps it as
# literal text instead of re-parsing it as a date.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$xlPasteValues = -4163
$scratch = $ws.Range("D1")

# Row 6
$scratch.Formula = '="2024-10-04"'
$scratch.Copy()
$ws.Range("A6").PasteSpecial($xlPasteValues)
$scratch.Clear()
$ws.Range("B6").Value = 0.003232

# Row 7
$scratch.Formula = '="2024-10-05"'
$scratch.Copy()
$ws.Range("A7").PasteSpecial($xlPasteValues)
$scratch.Clear()
$ws.Range("B7").Value = 0.003232
